$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Add a new worksheet right after "bf1" and rename it "Sheet1"
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet1"

# Column widths A:G (~22.71 characters wide in the source file)
$ws2.Range("A1:G10").ColumnWidth = 21.86

# A couple of rows got manually resized by the author
$ws2.Rows.Item(5).RowHeight = 15.75
$ws2.Rows.Item(6).RowHeight = 15

# Room / passage labels for the new map grid, entered in the order the
# author actually typed them (controls shared-string table order)
$entries = @(
    @(5, 7, "dungeon"),
    @(8, 1, "torture_chamber"),
    @(4, 6, "guard_room"),
    @(2, 1, "giant_centipede_room"),
    @(4, 3, "demon_bat_room"),
    @(3, 5, "skeleton_room"),
    @(3, 6, "empty_hallway"),
    @(4, 1, "large_rat_room"),
    @(8, 6, "pre_flooded_hall"),
    @(9, 6, "flooded_hall"),
    @(10, 6, "ritual_room"),
    @(1, 1, "stairs"),
    @(3, 1, "empty_passageway"),
    @(3, 3, "empty_room"),
    @(3, 4, "empty_passageway"),
    @(4, 2, "empty_passageway"),
    @(5, 3, "empty_passageway"),
    @(5, 5, "empty_passageway"),
    @(5, 6, "goblin_room"),
    @(6, 3, "empty_passageway"),
    @(6, 4, "empty_passageway"),
    @(6, 5, "giant_spider_room"),
    @(7, 3, "guard_room"),
    @(8, 2, "goblin_room"),
    @(8, 3, "empty_passageway"),
    @(8, 4, "empty_passageway"),
    @(8, 5, "giant_spider_room"),
    @(9, 3, "dungeon"),
    @(10, 3, "jail_cell")
)

foreach ($entry in $entries) {
    $ws2.Cells.Item($entry[0], $entry[1]).Value = $entry[2]
}

# Every cell in A1:G10 gets a medium box border and word-wrap
for ($r = 1; $r -le 10; $r++) {
    for ($c = 1; $c -le 7; $c++) {
        $cell = $ws2.Cells.Item($r, $c)
        $cell.Borders.Weight = -4138
        $cell.WrapText = $true
    }
}

# Author's last selection before saving
$ws2.Range("F11").Select()

# The new sheet is the active tab
$ws2.Activate()

Write-Output "done"
